$d = $word.ActiveDocument

function Replace-Needle($needle, $replacement) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($needle)
    if ($idx -lt 0) {
        throw "Needle not found: $needle"
    }
    $start = $idx
    $end = $idx + $needle.Length
    $rng = $d.Range($start, $end)
    $rng.Delete()
    $ins = $d.Range($start, $start)
    $ins.InsertBefore($replacement)
}

# 1) Title meta tag: "1 Introduction to Game Maker" -> "5 The Wall"
Replace-Needle "1 Introduction to Game Maker" "5 The Wall"

# 2) Description meta tag: collapse the three runs (highlighted "This ",
#    "is our first article...it", ". ") into new plain text.
Replace-Needle "This is our first article in a new series of Game Maker, where we will be introducing you to it. " "In this tutorial we will be looking at how we can create a wall for the room."

# 3) Revised meta tag: collapse "Thursday 30" + superscript "th" + ", 2025"
#    into "November 20, 2025".
Replace-Needle "Thursday 30th, 2025" "November 20, 2025"

# 4) URL meta tag: update the article path.
Replace-Needle "Enlightenment/Articles/2025/4_Game_Maker/1_Introduction_to_Game_Maker/1_Introduction_to_GameMaker.html" "Enlightenment/Articles/2025/4_Game_Maker/5_The_Wall/5_The_Wall.html"

"done"
